$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.206.53'
$ws.Range('E2').Value = '  -4.20%  '
$ws.Range('D3').Value = '3.304.36'
$ws.Range('E3').Value = '  -5.37%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.91'
$ws.Range('E5').Value = '  -4.72%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '611.31'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('E7').Value = '  -6.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.381'
$ws.Range('E8').Value = '  -6.41%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.945'
$ws.Range('E10').Value = '  -6.23%  '
$ws.Range('D11').Value = '3.301.03'
$ws.Range('E11').Value = '  -5.39%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '41.52'
$ws.Range('E12').Value = '  -3.74%  '
$ws.Range('E13').Value = '  -3.44%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.95'
$ws.Range('E14').Value = '  -4.45%  '
$ws.Range('D15').Value = '91.086.66'
$ws.Range('E15').Value = '  -4.09%  '
$ws.Range('D16').Value = '3.925.67'
$ws.Range('E16').Value = '  -5.42%  '
$ws.Range('E17').Value = '  -6.34%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.02'
$ws.Range('E18').Value = '  -6.10%  '
$ws.Range('D19').Value = '3.304.71'
$ws.Range('E19').Value = '  -5.62%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.09'
$ws.Range('E20').Value = '  -5.81%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.79'
$ws.Range('E21').Value = '  -6.38%  '
$ws.Range('E22').Value = '  +5.49%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '486.12'
$ws.Range('E23').Value = '  -4.34%  '
$ws.Range('E24').Value = '  -14.43%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0000179'
$ws.Range('E25').Value = '  -7.84%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.17'
$ws.Range('E26').Value = '  -8.23%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '88.92'
$ws.Range('E27').Value = '  -7.09%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.70'
$ws.Range('E28').Value = '  -4.52%  '
$ws.Range('D29').Value = '3.482.11'
$ws.Range('E29').Value = '  -5.32%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '10.90'
$ws.Range('E31').Value = '  -9.03%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.136'
$ws.Range('E32').Value = '  -1.77%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.60'
$ws.Range('E33').Value = '  -6.45%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.996'
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('E35').Value = '  -7.89%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '27.87'
$ws.Range('E36').Value = '  -10.44%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.519'
$ws.Range('E37').Value = '  -10.05%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '539.43'
$ws.Range('E38').Value = '  -3.49%  '
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '7.28'
$ws.Range('E40').Value = '  -6.89%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.35'
$ws.Range('E41').Value = '  -8.46%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.146'
$ws.Range('E42').Value = '  -3.62%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.855'
$ws.Range('E43').Value = '  -9.41%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '23.69'
$ws.Range('E44').Value = '  -1.60%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.60'
$ws.Range('E45').Value = '  +1.60%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0406'
$ws.Range('E46').Value = '  -3.21%  '
$ws.Range('B47').Value = 'ImmutableX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.64'
$ws.Range('E47').Value = '  -4.63%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.33'
$ws.Range('E48').Value = '  -6.98%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.08'
$ws.Range('E49').Value = '  -5.15%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.90'
$ws.Range('E50').Value = '  -2.85%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '51.42'
$ws.Range('E51').Value = '  -3.86%  '
